# Update "excel" / "excel_selected" / "duplicated" columns (P, Q, R) of the
# dataframe that feeds the bar chart on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("P4").Value = 15
$ws.Range("P5").Value = 21

$ws.Range("P7").Value = 65
$ws.Range("Q7").Value = 1

$ws.Range("P8").Value = 2
$ws.Range("Q8").Value = 1
$ws.Range("R8").Value = 77

$ws.Range("P9").Value = 466
$ws.Range("Q9").Value = 11
